$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "5+59=64"
$t.Cell(1,2).Range.Text = "94-55=39"
$t.Cell(1,3).Range.Text = "28+3=31"
$t.Cell(1,4).Range.Text = "25+8=33"
$t.Cell(1,5).Range.Text = "49+4=53"
$t.Cell(2,1).Range.Text = "57+18=75"
$t.Cell(2,2).Range.Text = "17+39=56"
$t.Cell(2,3).Range.Text = "32-3=29"
$t.Cell(2,4).Range.Text = "19+55=74"
$t.Cell(2,5).Range.Text = "17+76=93"
$t.Cell(3,1).Range.Text = "38+58=96"
$t.Cell(3,2).Range.Text = "57-18=39"
$t.Cell(3,3).Range.Text = "24-5=19"
$t.Cell(3,4).Range.Text = "15+7=22"
$t.Cell(3,5).Range.Text = "55-28=27"
$t.Cell(4,1).Range.Text = "15+78=93"
$t.Cell(4,2).Range.Text = "92-75=17"
$t.Cell(4,3).Range.Text = "80-26=54"
$t.Cell(4,4).Range.Text = "33-26=7"
$t.Cell(4,5).Range.Text = "87+7=94"
$t.Cell(5,1).Range.Text = "52-14=38"
$t.Cell(5,2).Range.Text = "81-16=65"
$t.Cell(5,3).Range.Text = "39+12=51"
$t.Cell(5,4).Range.Text = "76-19=57"
$t.Cell(5,5).Range.Text = "26+49=75"
$t.Cell(6,1).Range.Text = "72-48=24"
$t.Cell(6,2).Range.Text = "62-17=45"
$t.Cell(6,3).Range.Text = "50-5=45"
$t.Cell(6,4).Range.Text = "19+39=58"
$t.Cell(6,5).Range.Text = "86-9=77"
$t.Cell(7,1).Range.Text = "24+29=53"
$t.Cell(7,2).Range.Text = "86+5=91"
$t.Cell(7,3).Range.Text = "29+37=66"
$t.Cell(7,4).Range.Text = "83-26=57"
$t.Cell(7,5).Range.Text = "70-14=56"
$t.Cell(8,1).Range.Text = "69+15=84"
$t.Cell(8,2).Range.Text = "69+17=86"
$t.Cell(8,3).Range.Text = "12+39=51"
$t.Cell(8,4).Range.Text = "41-15=26"
$t.Cell(8,5).Range.Text = "89+9=98"
$t.Cell(9,1).Range.Text = "9+52=61"
$t.Cell(9,2).Range.Text = "57-39=18"
$t.Cell(9,3).Range.Text = "90-83=7"
$t.Cell(9,4).Range.Text = "79+7=86"
$t.Cell(9,5).Range.Text = "84-39=45"
$t.Cell(10,1).Range.Text = "83-55=28"
$t.Cell(10,2).Range.Text = "76-69=7"
$t.Cell(10,3).Range.Text = "62-38=24"
$t.Cell(10,4).Range.Text = "64+28=92"
$t.Cell(10,5).Range.Text = "3+58=61"
$t.Cell(11,1).Range.Text = "33-25=8"
$t.Cell(11,2).Range.Text = "39+54=93"
$t.Cell(11,3).Range.Text = "9+17=26"
$t.Cell(11,4).Range.Text = "94-19=75"
$t.Cell(11,5).Range.Text = "74-59=15"
$t.Cell(12,1).Range.Text = "83-57=26"
$t.Cell(12,2).Range.Text = "61-8=53"
$t.Cell(12,3).Range.Text = "17+9=26"
$t.Cell(12,4).Range.Text = "34-18=16"
$t.Cell(12,5).Range.Text = "29+4=33"
$t.Cell(13,1).Range.Text = "96-49=47"
$t.Cell(13,2).Range.Text = "28+59=87"
$t.Cell(13,3).Range.Text = "20-6=14"
$t.Cell(13,4).Range.Text = "79+6=85"
$t.Cell(13,5).Range.Text = "81-12=69"
$t.Cell(14,1).Range.Text = "70-62=8"
$t.Cell(14,2).Range.Text = "30-13=17"
$t.Cell(14,3).Range.Text = "20-6=14"
$t.Cell(14,4).Range.Text = "70-25=45"
$t.Cell(14,5).Range.Text = "69+5=74"
$t.Cell(15,1).Range.Text = "55-37=18"
$t.Cell(15,2).Range.Text = "29+4=33"
$t.Cell(15,3).Range.Text = "59+37=96"
$t.Cell(15,4).Range.Text = "40-28=12"
$t.Cell(15,5).Range.Text = "68-29=39"
$t.Cell(16,1).Range.Text = "3+59=62"
$t.Cell(16,2).Range.Text = "55+17=72"
$t.Cell(16,3).Range.Text = "72-55=17"
$t.Cell(16,4).Range.Text = "88-19=69"
$t.Cell(16,5).Range.Text = "19+73=92"
$t.Cell(17,1).Range.Text = "17+66=83"
$t.Cell(17,2).Range.Text = "95-17=78"
$t.Cell(17,3).Range.Text = "28+58=86"
$t.Cell(17,4).Range.Text = "7+85=92"
$t.Cell(17,5).Range.Text = "26+6=32"
$t.Cell(18,1).Range.Text = "41-2=39"
$t.Cell(18,2).Range.Text = "94-66=28"
$t.Cell(18,3).Range.Text = "60-22=38"
$t.Cell(18,4).Range.Text = "95-49=46"
$t.Cell(18,5).Range.Text = "60-13=47"
$t.Cell(19,1).Range.Text = "6+75=81"
$t.Cell(19,2).Range.Text = "56+15=71"
$t.Cell(19,3).Range.Text = "60-14=46"
$t.Cell(19,4).Range.Text = "60-9=51"
$t.Cell(19,5).Range.Text = "17+47=64"
$t.Cell(20,1).Range.Text = "27+69=96"
$t.Cell(20,2).Range.Text = "22+59=81"
$t.Cell(20,3).Range.Text = "28+63=91"
$t.Cell(20,4).Range.Text = "6+36=42"
$t.Cell(20,5).Range.Text = "32-5=27"
